$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths
# (the COM ColumnWidth setter adds a constant 5/6-character padding before it is
# persisted as the OOXML "width" attribute, so we compensate by subtracting it
# here to land exactly on the target stored widths of 14 and 33)
$ws.Columns.Item(1).ColumnWidth = 14 - 0.8333333333333334
$ws.Columns.Item(3).ColumnWidth = 33 - 0.8333333333333334

# Update existing row 2 data
$ws.Cells.Item(2, 1).Value = "GR3922"
$ws.Cells.Item(2, 2).Value = "sku invalido"
$ws.Cells.Item(2, 3).Value = "produccion-lerma-semana-05.xlsx"

# Data for new rows 3-13
$data = @(
    @("GR3890", "sku invalido", "produccion-lerma-semana-05.xlsx"),
    @("ANSAGR4051", "sku invalido", "produccion-lerma-semana-05.xlsx"),
    @("ANSAH340", "sku invalido", "produccion-lerma-semana-05.xlsx"),
    @(10481573, "sku invalido", "produccion-lerma-semana-05.xlsx"),
    @("GIFT81747334", "sku invalido", "produccion-lerma-semana-05.xlsx"),
    @("GIFT80863470", "sku invalido", "produccion-lerma-semana-05.xlsx"),
    @(11081214, "sku invalido", "produccion-lerma-semana-05.xlsx"),
    @(10595873, "sku invalido", "produccion-lerma-semana-05.xlsx"),
    @("GIFT79485824", "sku invalido", "produccion-lerma-semana-05.xlsx"),
    @(10595875, "sku invalido", "produccion-lerma-semana-05.xlsx"),
    @("GIFT12659082", "sku invalido", "produccion-lerma-semana-05.xlsx")
)

$row = 3
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
